$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LIQ_dict")

# --- Row 2: title ---
$ws.Range("B2").Value = "Fragebogen: Literarische Interessen"

# --- Row 5-8: TLIQ_0003 prompt/options (reading/listening -> also watching for plays) ---
$ws.Range("B5").Value = "Wie gern lesen oder hören (z.B. Audiobuch, Lesung) Sie diese literarischen Gattungen?"
$ws.Range("B6").Value = "Wie gern lesen oder hören (z.B. Audiobuch, Lesung) Sie literarische **Prosa** (z.B. Romane, Erzählungen, Kurzgeschichten)?"
$ws.Range("B7").Value = "Wie gern lesen oder hören (z.B. Audiobuch, Lesung) Sie **Gedichte**?"
$ws.Range("B8").Value = "Wie gern lesen oder sehen (z.B. Theatervorstellung) Sie **Dramen** (z.B. Komödien, Tragödien)?"
$ws.Range("C8").Value = "How much do you like reading or watching (e.g. theatre performance) to plays (e.g., comedies, tragedies)?"

# --- Row 24-25: TLIQ_0005 genre options (reading peak) ---
$ws.Range("B24").Value = "Literarische Prosa (z.B. Romane, Erzählungen, Kurzgeschichten)"
$ws.Range("C24").Value = "Literary Prose (e.g. novels, short stories)"
$ws.Range("B25").Value = "Gedichte"
$ws.Range("C25").Value = "Poems"

# --- Row 32-33: TLIQ_0006 genre options (writing) ---
$ws.Range("B32").Value = "Literarische Prosa (z.B. Romane, Erzählungen, Kurzgeschichten)"
$ws.Range("B33").Value = "Gedichte"
$ws.Range("C33").Value = "Poems"

# --- Row 35: TLIQ_0007 prompt (peak-of-interest reading) ---
$ws.Range("B35").Value = "Denken Sie nun bitte einmal an **die Zeit in Ihrem Leben, in der Sie sich am meisten mit Literatur beschäftigt haben**. \\ Wie viele Stunden haben Sie auf dem Höhepunkt Ihres literarischen Interesses durchschnittlich in einer Woche damit verbracht, die folgenden literarischen Gattungen zu **lesen**."

# --- Row 37-38: TLIQ_0007 genre options (peak-of-interest reading) ---
$ws.Range("B37").Value = "Literarische Prosa (z.B. Romane, Erzählungen, Kurzgeschichten)"
$ws.Range("B38").Value = "Gedichte"
$ws.Range("C38").Value = "Poems"

# --- Row 42-43: TLIQ_0008 genre options (peak-of-interest writing) ---
$ws.Range("B42").Value = "Literarische Prosa (z.B. Romane, Erzählungen, Kurzgeschichten)"
$ws.Range("B43").Value = "Gedichte"

# --- Selection change ---
$ws.Range("B42:B43").Select()
